$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.580195426940918
$ws.Range("B1").Value = 2.315661191940308
$ws.Range("C1").Value = 4.517023086547852
$ws.Range("D1").Value = 4.729501247406006
$ws.Range("E1").Value = 1.567581295967102
